$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -22.305
$ws.Range("A10").Value = -21.777
$ws.Range("A12").Value = -21.546
$ws.Range("A18").Value = -22.095
$ws.Range("A25").Value = -21.765
$ws.Range("A37").Value = -20.217
$ws.Range("A55").Value = -22.184
$ws.Range("A68").Value = -21.58300000000001
$ws.Range("A77").Value = -20.475
$ws.Range("A78").Value = -20.095
$ws.Range("A79").Value = -21.246
$ws.Range("A80").Value = -20.239
$ws.Range("A81").Value = -21.78
$ws.Range("A82").Value = -22.067
$ws.Range("A84").Value = -21.988
$ws.Range("A101").Value = -20.672
$ws.Range("A102").Value = -20.353
